# Update LR-pair rows for Efnb2-Epha3 with refreshed TPM-based NATMI output.
# Source table grows from 3 sender clusters x 4 target clusters (12 rows)
# to 4 sender clusters x 4 target clusters (16 rows), adding 'Resolving-Mac'
# as a sending cluster (rows 14-17) and refreshing all previously computed
# ligand/receptor/edge statistics in rows 2-13 with newly recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Efnb2"
$ws.Cells.Item(2,3).Value = "Epha3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 45.71598933333333
$ws.Cells.Item(2,8).Value = 137.147968
$ws.Cells.Item(2,9).Value = 0.6549002937372808
$ws.Cells.Item(2,10).Value = 0.6549002937372808
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.009849666666666666
$ws.Cells.Item(2,14).Value = 0.029549
$ws.Cells.Item(2,15).Value = 0.0002013876315934659
$ws.Cells.Item(2,16).Value = 0.0002013876315934659
$ws.Cells.Item(2,17).Value = 0.4502872562702222
$ws.Cells.Item(2,18).Value = 4.052585306431999
$ws.Cells.Item(2,19).Value = 0.0001318888190856161
$ws.Cells.Item(2,20).Value = 0.0001318888190856161

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Efnb2"
$ws.Cells.Item(3,3).Value = "Epha3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 45.71598933333333
$ws.Cells.Item(3,8).Value = 137.147968
$ws.Cells.Item(3,9).Value = 0.6549002937372808
$ws.Cells.Item(3,10).Value = 0.6549002937372808
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 46.25093466666667
$ws.Cells.Item(3,14).Value = 138.752804
$ws.Cells.Item(3,15).Value = 0.9456529349389956
$ws.Cells.Item(3,16).Value = 0.9456529349389956
$ws.Cells.Item(3,17).Value = 2114.40723587803
$ws.Cells.Item(3,18).Value = 19029.66512290227
$ws.Cells.Item(3,19).Value = 0.6193083848650699
$ws.Cells.Item(3,20).Value = 0.6193083848650699

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Efnb2"
$ws.Cells.Item(4,3).Value = "Epha3"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 45.71598933333333
$ws.Cells.Item(4,8).Value = 137.147968
$ws.Cells.Item(4,9).Value = 0.6549002937372808
$ws.Cells.Item(4,10).Value = 0.6549002937372808
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.607896333333333
$ws.Cells.Item(4,14).Value = 7.823689
$ws.Cells.Item(4,15).Value = 0.05332140505715427
$ws.Cells.Item(4,16).Value = 0.05332140505715428
$ws.Cells.Item(4,17).Value = 119.2225609571058
$ws.Cells.Item(4,18).Value = 1073.003048613952
$ws.Cells.Item(4,19).Value = 0.03492020383441486
$ws.Cells.Item(4,20).Value = 0.03492020383441486

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Efnb2"
$ws.Cells.Item(5,3).Value = "Epha3"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 45.71598933333333
$ws.Cells.Item(5,8).Value = 137.147968
$ws.Cells.Item(5,9).Value = 0.6549002937372808
$ws.Cells.Item(5,10).Value = 0.6549002937372808
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.04031433333333333
$ws.Cells.Item(5,14).Value = 0.120943
$ws.Cells.Item(5,15).Value = 0.0008242723722565416
$ws.Cells.Item(5,16).Value = 0.0008242723722565415
$ws.Cells.Item(5,17).Value = 1.843009632647111
$ws.Cells.Item(5,18).Value = 16.587086693824
$ws.Cells.Item(5,19).Value = 0.0005398162187103343
$ws.Cells.Item(5,20).Value = 0.0005398162187103343

# Row 6: FAPs -> ECs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Efnb2"
$ws.Cells.Item(6,3).Value = "Epha3"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 12.691493
$ws.Cells.Item(6,8).Value = 38.074479
$ws.Cells.Item(6,9).Value = 0.1818108415648851
$ws.Cells.Item(6,10).Value = 0.1818108415648851
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.009849666666666666
$ws.Cells.Item(6,14).Value = 0.029549
$ws.Cells.Item(6,15).Value = 0.0002013876315934659
$ws.Cells.Item(6,16).Value = 0.0002013876315934659
$ws.Cells.Item(6,17).Value = 0.1250069755523333
$ws.Cells.Item(6,18).Value = 1.125062779971
$ws.Cells.Item(6,19).Value = 0.00003661445478076707
$ws.Cells.Item(6,20).Value = 0.00003661445478076707

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Efnb2"
$ws.Cells.Item(7,3).Value = "Epha3"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 12.691493
$ws.Cells.Item(7,8).Value = 38.074479
$ws.Cells.Item(7,9).Value = 0.1818108415648851
$ws.Cells.Item(7,10).Value = 0.1818108415648851
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 46.25093466666667
$ws.Cells.Item(7,14).Value = 138.752804
$ws.Cells.Item(7,15).Value = 0.9456529349389956
$ws.Cells.Item(7,16).Value = 0.9456529349389956
$ws.Cells.Item(7,17).Value = 586.9934135654573
$ws.Cells.Item(7,18).Value = 5282.940722089115
$ws.Cells.Item(7,19).Value = 0.1719299559295623
$ws.Cells.Item(7,20).Value = 0.1719299559295623

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Efnb2"
$ws.Cells.Item(8,3).Value = "Epha3"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 12.691493
$ws.Cells.Item(8,8).Value = 38.074479
$ws.Cells.Item(8,9).Value = 0.1818108415648851
$ws.Cells.Item(8,10).Value = 0.1818108415648851
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.607896333333333
$ws.Cells.Item(8,14).Value = 7.823689
$ws.Cells.Item(8,15).Value = 0.05332140505715427
$ws.Cells.Item(8,16).Value = 0.05332140505715428
$ws.Cells.Item(8,17).Value = 33.09809805922566
$ws.Cells.Item(8,18).Value = 297.882882533031
$ws.Cells.Item(8,19).Value = 0.009694409526863337
$ws.Cells.Item(8,20).Value = 0.009694409526863337

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Efnb2"
$ws.Cells.Item(9,3).Value = "Epha3"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 12.691493
$ws.Cells.Item(9,8).Value = 38.074479
$ws.Cells.Item(9,9).Value = 0.1818108415648851
$ws.Cells.Item(9,10).Value = 0.1818108415648851
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.04031433333333333
$ws.Cells.Item(9,14).Value = 0.120943
$ws.Cells.Item(9,15).Value = 0.0008242723722565416
$ws.Cells.Item(9,16).Value = 0.0008242723722565415
$ws.Cells.Item(9,17).Value = 0.5116490792996666
$ws.Cells.Item(9,18).Value = 4.604841713697
$ws.Cells.Item(9,19).Value = 0.0001498616536786461
$ws.Cells.Item(9,20).Value = 0.000149861653678646

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Efnb2"
$ws.Cells.Item(10,3).Value = "Epha3"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 11.24784666666667
$ws.Cells.Item(10,8).Value = 33.74354
$ws.Cells.Item(10,9).Value = 0.161130015850732
$ws.Cells.Item(10,10).Value = 0.161130015850732
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.009849666666666666
$ws.Cells.Item(10,14).Value = 0.029549
$ws.Cells.Item(10,15).Value = 0.0002013876315934659
$ws.Cells.Item(10,16).Value = 0.0002013876315934659
$ws.Cells.Item(10,17).Value = 0.1107875403844444
$ws.Cells.Item(10,18).Value = 0.9970878634599999
$ws.Cells.Item(10,19).Value = 0.00003244959227079653
$ws.Cells.Item(10,20).Value = 0.00003244959227079653

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Efnb2"
$ws.Cells.Item(11,3).Value = "Epha3"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 11.24784666666667
$ws.Cells.Item(11,8).Value = 33.74354
$ws.Cells.Item(11,9).Value = 0.161130015850732
$ws.Cells.Item(11,10).Value = 0.161130015850732
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 46.25093466666667
$ws.Cells.Item(11,14).Value = 138.752804
$ws.Cells.Item(11,15).Value = 0.9456529349389956
$ws.Cells.Item(11,16).Value = 0.9456529349389956
$ws.Cells.Item(11,17).Value = 520.2234213206845
$ws.Cells.Item(11,18).Value = 4682.010791886159
$ws.Cells.Item(11,19).Value = 0.1523730723960116
$ws.Cells.Item(11,20).Value = 0.1523730723960116

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12,1).Value = "MuSCs"
$ws.Cells.Item(12,2).Value = "Efnb2"
$ws.Cells.Item(12,3).Value = "Epha3"
$ws.Cells.Item(12,4).Value = "MuSCs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 11.24784666666667
$ws.Cells.Item(12,8).Value = 33.74354
$ws.Cells.Item(12,9).Value = 0.161130015850732
$ws.Cells.Item(12,10).Value = 0.161130015850732
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.607896333333333
$ws.Cells.Item(12,14).Value = 7.823689
$ws.Cells.Item(12,15).Value = 0.05332140505715427
$ws.Cells.Item(12,16).Value = 0.05332140505715428
$ws.Cells.Item(12,17).Value = 29.33321807989555
$ws.Cells.Item(12,18).Value = 263.99896271906
$ws.Cells.Item(12,19).Value = 0.008591678842042568
$ws.Cells.Item(12,20).Value = 0.008591678842042568

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13,1).Value = "MuSCs"
$ws.Cells.Item(13,2).Value = "Efnb2"
$ws.Cells.Item(13,3).Value = "Epha3"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 11.24784666666667
$ws.Cells.Item(13,8).Value = 33.74354
$ws.Cells.Item(13,9).Value = 0.161130015850732
$ws.Cells.Item(13,10).Value = 0.161130015850732
$ws.Cells.Item(13,11).Value = 1
$ws.Cells.Item(13,12).Value = 0.3333333333333333
$ws.Cells.Item(13,13).Value = 0.04031433333333333
$ws.Cells.Item(13,14).Value = 0.120943
$ws.Cells.Item(13,15).Value = 0.0008242723722565416
$ws.Cells.Item(13,16).Value = 0.0008242723722565415
$ws.Cells.Item(13,17).Value = 0.4534494398022222
$ws.Cells.Item(13,18).Value = 4.08104495822
$ws.Cells.Item(13,19).Value = 0.000132815020407017
$ws.Cells.Item(13,20).Value = 0.000132815020407017

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14,1).Value = "Resolving-Mac"
$ws.Cells.Item(14,2).Value = "Efnb2"
$ws.Cells.Item(14,3).Value = "Epha3"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 0.6666666666666666
$ws.Cells.Item(14,7).Value = 0.1507006666666667
$ws.Cells.Item(14,8).Value = 0.452102
$ws.Cells.Item(14,9).Value = 0.00215884884710222
$ws.Cells.Item(14,10).Value = 0.00215884884710222
$ws.Cells.Item(14,11).Value = 1
$ws.Cells.Item(14,12).Value = 0.3333333333333333
$ws.Cells.Item(14,13).Value = 0.009849666666666666
$ws.Cells.Item(14,14).Value = 0.029549
$ws.Cells.Item(14,15).Value = 0.0002013876315934659
$ws.Cells.Item(14,16).Value = 0.0002013876315934659
$ws.Cells.Item(14,17).Value = 0.001484351333111111
$ws.Cells.Item(14,18).Value = 0.013359161998
$ws.Cells.Item(14,19).Value = 0.0000004347654562862005
$ws.Cells.Item(14,20).Value = 0.0000004347654562862004

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15,1).Value = "Resolving-Mac"
$ws.Cells.Item(15,2).Value = "Efnb2"
$ws.Cells.Item(15,3).Value = "Epha3"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 0.6666666666666666
$ws.Cells.Item(15,7).Value = 0.1507006666666667
$ws.Cells.Item(15,8).Value = 0.452102
$ws.Cells.Item(15,9).Value = 0.00215884884710222
$ws.Cells.Item(15,10).Value = 0.00215884884710222
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 46.25093466666667
$ws.Cells.Item(15,14).Value = 138.752804
$ws.Cells.Item(15,15).Value = 0.9456529349389956
$ws.Cells.Item(15,16).Value = 0.9456529349389956
$ws.Cells.Item(15,17).Value = 6.970046688223111
$ws.Cells.Item(15,18).Value = 62.730420194008
$ws.Cells.Item(15,19).Value = 0.002041521748351881
$ws.Cells.Item(15,20).Value = 0.002041521748351881

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16,1).Value = "Resolving-Mac"
$ws.Cells.Item(16,2).Value = "Efnb2"
$ws.Cells.Item(16,3).Value = "Epha3"
$ws.Cells.Item(16,4).Value = "MuSCs"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 0.6666666666666666
$ws.Cells.Item(16,7).Value = 0.1507006666666667
$ws.Cells.Item(16,8).Value = 0.452102
$ws.Cells.Item(16,9).Value = 0.00215884884710222
$ws.Cells.Item(16,10).Value = 0.00215884884710222
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 2.607896333333333
$ws.Cells.Item(16,14).Value = 7.823689
$ws.Cells.Item(16,15).Value = 0.05332140505715427
$ws.Cells.Item(16,16).Value = 0.05332140505715428
$ws.Cells.Item(16,17).Value = 0.3930117160308889
$ws.Cells.Item(16,18).Value = 3.537105444278
$ws.Cells.Item(16,19).Value = 0.000115112853833508
$ws.Cells.Item(16,20).Value = 0.000115112853833508

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17,1).Value = "Resolving-Mac"
$ws.Cells.Item(17,2).Value = "Efnb2"
$ws.Cells.Item(17,3).Value = "Epha3"
$ws.Cells.Item(17,4).Value = "Resolving-Mac"
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 0.6666666666666666
$ws.Cells.Item(17,7).Value = 0.1507006666666667
$ws.Cells.Item(17,8).Value = 0.452102
$ws.Cells.Item(17,9).Value = 0.00215884884710222
$ws.Cells.Item(17,10).Value = 0.00215884884710222
$ws.Cells.Item(17,11).Value = 1
$ws.Cells.Item(17,12).Value = 0.3333333333333333
$ws.Cells.Item(17,13).Value = 0.04031433333333333
$ws.Cells.Item(17,14).Value = 0.120943
$ws.Cells.Item(17,15).Value = 0.0008242723722565416
$ws.Cells.Item(17,16).Value = 0.0008242723722565415
$ws.Cells.Item(17,17).Value = 0.006075396909555556
$ws.Cells.Item(17,18).Value = 0.05467857218599999
$ws.Cells.Item(17,19).Value = 0.000001779479460544247
$ws.Cells.Item(17,20).Value = 0.000001779479460544246
